$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices in column D are stored as plain text (e.g. "29.086.92", "1.140")
# even though several of them would otherwise be auto-parsed as numbers by
# Excel (losing trailing zeros, e.g. "1.140" -> 1.14). Force text entry by
# switching the cell to the Text number format before assigning the value,
# then restore the default "Normal" style so no visible formatting changes.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.086.92"
$ws.Range("E2").Value = "  +0.19%  "
Set-TextValue $ws.Range("D3") "1.837.00"
$ws.Range("E3").Value = "  +0.46%  "
Set-TextValue $ws.Range("D4") "0.9996"
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws.Range("D5") "243.39"
$ws.Range("E5").Value = "  +0.72%  "
Set-TextValue $ws.Range("D6") "0.6277"
$ws.Range("E6").Value = "  -1.37%  "
Set-TextValue $ws.Range("D7") "1.002"
$ws.Range("E7").Value = "  +0.14%  "
Set-TextValue $ws.Range("D8") "0.07602"
$ws.Range("E8").Value = "  +3.64%  "
Set-TextValue $ws.Range("D9") "0.2935"
$ws.Range("E9").Value = "  +0.19%  "
Set-TextValue $ws.Range("D10") "22.62"
$ws.Range("E10").Value = "  -1.27%  "
Set-TextValue $ws.Range("D11") "0.07747"
$ws.Range("E11").Value = "  +1.05%  "
Set-TextValue $ws.Range("D12") "1.846.17"
$ws.Range("E12").Value = "  +0.99%  "
Set-TextValue $ws.Range("D13") "4.966"
$ws.Range("E13").Value = "  -0.47%  "
Set-TextValue $ws.Range("D14") "0.6654"
$ws.Range("E14").Value = "  +0.05%  "
Set-TextValue $ws.Range("D15") "0.00001009"
$ws.Range("E15").Value = "  +17.22%  "
Set-TextValue $ws.Range("D16") "82.89"
$ws.Range("E16").Value = "  +0.80%  "
Set-TextValue $ws.Range("D17") "6.065"
$ws.Range("E17").Value = "  -0.13%  "
Set-TextValue $ws.Range("D18") "29.101.03"
$ws.Range("E18").Value = "  +0.29%  "
Set-TextValue $ws.Range("D19") "227.33"
$ws.Range("E19").Value = "  +1.48%  "
Set-TextValue $ws.Range("D20") "12.39"
$ws.Range("E20").Value = "  -0.13%  "
Set-TextValue $ws.Range("D21") "1.003"
$ws.Range("E21").Value = "  +0.29%  "
Set-TextValue $ws.Range("D22") "7.211"
$ws.Range("E22").Value = "  +1.56%  "
Set-TextValue $ws.Range("D23") "1.001"
$ws.Range("E23").Value = "  +0.13%  "
Set-TextValue $ws.Range("D24") "159.31"
$ws.Range("E24").Value = "  +0.75%  "
Set-TextValue $ws.Range("D25") "8.525"
$ws.Range("E25").Value = "  +0.74%  "
Set-TextValue $ws.Range("D26") "0.1385"
$ws.Range("E26").Value = "  +0.55%  "
Set-TextValue $ws.Range("D27") "17.96"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.44%  "
Set-TextValue $ws.Range("D30") "4.025"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  -0.33%  "
Set-TextValue $ws.Range("D32") "0.05252"
$ws.Range("E32").Value = "  -1.03%  "
Set-TextValue $ws.Range("D33") "1.845"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  -1.01%  "
Set-TextValue $ws.Range("D35") "1.140"
$ws.Range("E35").Value = "  -1.13%  "
Set-TextValue $ws.Range("D36") "2.701"
$ws.Range("E36").Value = "  +2.09%  "
Set-TextValue $ws.Range("D37") "1.245.83"
$ws.Range("E37").Value = "  -3.33%  "
Set-TextValue $ws.Range("D38") "2.765"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  +0.32%  "
Set-TextValue $ws.Range("D40") "6.364"
$ws.Range("E40").Value = "  +0.41%  "
Set-TextValue $ws.Range("D41") "0.8973"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E42").Value = "  +0.29%  "
Set-TextValue $ws.Range("D43") "102.23"
$ws.Range("E43").Value = "  -0.60%  "
Set-TextValue $ws.Range("D44") "1.976.20"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  -5.39%  "
Set-TextValue $ws.Range("D46") "64.44"
$ws.Range("E46").Value = "  +0.32%  "
Set-TextValue $ws.Range("D47") "0.5111"
$ws.Range("E47").Value = "  -0.47%  "
Set-TextValue $ws.Range("D48") "0.4049"
$ws.Range("E48").Value = "  +1.66%  "
Set-TextValue $ws.Range("D49") "8.860"
$ws.Range("E49").Value = "  +2.13%  "
Set-TextValue $ws.Range("D50") "0.05757"
$ws.Range("E50").Value = "  -1.23%  "
Set-TextValue $ws.Range("D51") "6.690"
$ws.Range("E51").Value = "  +0.10%  "
